$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -7.881
$ws.Range("E3").Value = 16.794
$ws.Range("B12").Value = 5.220999999999999
$ws.Range("D14").Value = -7.318000000000001
$ws.Range("E20").Value = 16.485
$ws.Range("E25").Value = 17.006
$ws.Range("D26").Value = -7.640000000000001
$ws.Range("B27").Value = 6.052
$ws.Range("E30").Value = 16.029
$ws.Range("D31").Value = -8.430000000000001
$ws.Range("B32").Value = 6.387
$ws.Range("D35").Value = -7.672
$ws.Range("B36").Value = 8.548
$ws.Range("D37").Value = -7.741
$ws.Range("B38").Value = 5.565
$ws.Range("E44").Value = 16.64
$ws.Range("D45").Value = -7.498
$ws.Range("B46").Value = 6.452
$ws.Range("E47").Value = 16.058
$ws.Range("D52").Value = -7.961
$ws.Range("B54").Value = 5.050000000000001
$ws.Range("B55").Value = 4.953
$ws.Range("B56").Value = 4.877000000000001
$ws.Range("D57").Value = -8.100999999999999
$ws.Range("E58").Value = 16.425
$ws.Range("B67").Value = 5.078
$ws.Range("B69").Value = 5.078
$ws.Range("B72").Value = 5.073
$ws.Range("E78").Value = 16.391
$ws.Range("D81").Value = -6.936
$ws.Range("B83").Value = 5.737
$ws.Range("D83").Value = -8.572000000000001
$ws.Range("E84").Value = 16.178
$ws.Range("B86").Value = 5.145999999999999
$ws.Range("E89").Value = 17.362
$ws.Range("B91").Value = 5.276
$ws.Range("E91").Value = 16.916
$ws.Range("E92").Value = 16.931
$ws.Range("B93").Value = 5.427000000000001
$ws.Range("E96").Value = 16.303
$ws.Range("B99").Value = 5.217000000000001
$ws.Range("D100").Value = -8.313000000000001
$ws.Range("D102").Value = -7.794
$ws.Range("E102").Value = 16.498
